$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Problem Solving" progress entry as row 17 of the M:Q block ---
# Seed row 17's formatting from row 14, which already carries the "closing"
# (top+bottom) table border that the new last row of Table4 needs - this mirrors
# how Excel re-closes the bottom border of a table when a row is appended.
$ws.Range("M14:Q14").Copy() | Out-Null
$ws.Range("M17:Q17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 14's date column used a slightly different (but equivalent) number format
# than the most recent rows; line row 17 up with row 16's date format instead.
$ws.Range("N17").NumberFormat = $ws.Range("N16").NumberFormat

# Now fill in the actual values for the new row.
$ws.Range("M17").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("N17").Value = 45102
$ws.Range("O17").Value = "1112.97/2200"
$ws.Range("P17").Value = 134464
$ws.Range("Q17").Formula = "=IF(ROW()>2,(`$P`$2-P17)/`$P`$2,""NA"")"

# --- Grow Table4 (M1:Q16) so it includes the newly-added row ---
$table4 = $ws.ListObjects.Item("Table4")
$table4.Resize($ws.Range("M1:Q17")) | Out-Null

# --- Match the saved selection location from the edited workbook ---
$ws.Range("O20").Select() | Out-Null
